$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C2:C7) from serial 45207 (2023-10-08)
# to serial 45208 (2023-10-09), preserving existing number formatting.
foreach ($row in 2..7) {
    $ws.Cells.Item($row, 3).Value2 = 45208
}
